# Auto-generated Excel COM-interop script to apply market-price/profit refresh
# (scheduled runner update) to the Cactuar_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1249.1765
$ws.Range("I80").Value = 834.7143
$ws.Range("J80").Value = 1539.3
$ws.Range("K80").Value = 2504.1429
$ws.Range("L80").Value = 4617.9
$ws.Range("M80").Value = -1506.1429
$ws.Range("N80").Value = -6613.9

$ws.Range("H83").Value = 1249.1765
$ws.Range("I83").Value = 834.7143
$ws.Range("J83").Value = 1539.3
$ws.Range("K83").Value = 7512.428699999999
$ws.Range("L83").Value = 13853.7
$ws.Range("M83").Value = -2520.428699999999
$ws.Range("N83").Value = -23837.7

$ws.Range("H86").Value = 1179078.4
$ws.Range("I86").Value = 1622281.2
$ws.Range("K86").Value = 1622281.2
$ws.Range("M86").Value = -1621158.2

$ws.Range("H89").Value = 1179078.4
$ws.Range("I89").Value = 1622281.2
$ws.Range("K89").Value = 8111406
$ws.Range("M89").Value = -8105790

$ws.Range("H100").Value = 1330
$ws.Range("I100").Value = 1330
$ws.Range("K100").Value = 1330
$ws.Range("M100").Value = -789

$ws.Range("H107").Value = 501.92856
$ws.Range("I107").Value = 321.85715
$ws.Range("J107").Value = 682
$ws.Range("K107").Value = 321.85715
$ws.Range("L107").Value = 682
$ws.Range("M107").Value = 1598.14285
$ws.Range("N107").Value = -4522

$ws.Range("H137").Value = 10101608
$ws.Range("I137").Value = 626372.8
$ws.Range("K137").Value = 1879118.4
$ws.Range("M137").Value = -1876568.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19798.51
$ws.Range("I32").Value = 21659.361
$ws.Range("J32").Value = 14645.385
$ws.Range("K32").Value = 21659.361
$ws.Range("L32").Value = 14645.385
$ws.Range("M32").Value = -21372.361
$ws.Range("N32").Value = -15219.385

$ws.Range("H45").Value = 2926.3333
$ws.Range("I45").Value = 1520.1666
$ws.Range("K45").Value = 1520.1666
$ws.Range("M45").Value = -1143.1666

$ws.Range("H61").Value = 10920.2
$ws.Range("I61").Value = 11876.765
$ws.Range("J61").Value = 5499.6665
$ws.Range("K61").Value = 11876.765
$ws.Range("L61").Value = 5499.6665
$ws.Range("M61").Value = -11664.765
$ws.Range("N61").Value = -5923.6665

$ws.Range("H74").Value = 7144228.5
$ws.Range("I74").Value = 9616448
$ws.Range("J74").Value = 2261.889
$ws.Range("K74").Value = 9616448
$ws.Range("L74").Value = 2261.889
$ws.Range("M74").Value = -9615574
$ws.Range("N74").Value = -4009.889

$ws.Range("H77").Value = 7144228.5
$ws.Range("I77").Value = 9616448
$ws.Range("J77").Value = 2261.889
$ws.Range("K77").Value = 48082240
$ws.Range("L77").Value = 11309.445
$ws.Range("M77").Value = -48077872
$ws.Range("N77").Value = -20045.445

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 12378.491
$ws.Range("I132").Value = 13331.234
$ws.Range("J132").Value = 6781.125
$ws.Range("K132").Value = 39993.702
$ws.Range("L132").Value = 20343.375
$ws.Range("M132").Value = -37463.702
$ws.Range("N132").Value = -25403.375

$ws.Range("H136").Value = 10920.2
$ws.Range("I136").Value = 11876.765
$ws.Range("J136").Value = 5499.6665
$ws.Range("K136").Value = 35630.295
$ws.Range("L136").Value = 16498.9995
$ws.Range("M136").Value = -33080.295
$ws.Range("N136").Value = -21598.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6869.729
$ws.Range("I31").Value = 2997
$ws.Range("K31").Value = 2997
$ws.Range("M31").Value = -2702

$ws.Range("H34").Value = 6869.729
$ws.Range("I34").Value = 2997
$ws.Range("K34").Value = 2997
$ws.Range("M34").Value = -2795

$ws.Range("H99").Value = 6042
$ws.Range("I99").Value = 4103.1816
$ws.Range("J99").Value = 7980.8184
$ws.Range("K99").Value = 4103.1816
$ws.Range("L99").Value = 7980.8184
$ws.Range("M99").Value = -2605.1816
$ws.Range("N99").Value = -10976.8184

$ws.Range("H107").Value = 1299163.2
$ws.Range("I107").Value = 1653252
$ws.Range("K107").Value = 1653252
$ws.Range("M107").Value = -1651332

$ws.Range("H126").Value = 6042
$ws.Range("I126").Value = 4103.1816
$ws.Range("J126").Value = 7980.8184
$ws.Range("K126").Value = 12309.5448
$ws.Range("L126").Value = 23942.4552
$ws.Range("M126").Value = -9839.5448
$ws.Range("N126").Value = -28882.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1383.5714
$ws.Range("J2").Value = 2649.6155
$ws.Range("L2").Value = 15897.693
$ws.Range("N2").Value = -16123.693

$ws.Range("H92").Value = 1355.8334
$ws.Range("I92").Value = 1348.5
$ws.Range("K92").Value = 4045.5
$ws.Range("M92").Value = -2797.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 79998.25
$ws.Range("J130").Value = 79998.25
$ws.Range("L130").Value = 79998.25
$ws.Range("N130").Value = -90038.25

$ws.Range("H132").Value = 403969.78
$ws.Range("I132").Value = 113633.22
$ws.Range("K132").Value = 340899.66
$ws.Range("M132").Value = -338369.66

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5373.25
$ws.Range("I61").Value = 4499
$ws.Range("J61").Value = 6247.5
$ws.Range("K61").Value = 4499
$ws.Range("L61").Value = 6247.5
$ws.Range("M61").Value = -4297
$ws.Range("N61").Value = -6651.5

$ws.Range("H93").Value = 1419.5
$ws.Range("I93").Value = 1416.7273
$ws.Range("K93").Value = 1416.7273
$ws.Range("M93").Value = -168.7273

$ws.Range("H113").Value = 5373.25
$ws.Range("I113").Value = 4499
$ws.Range("J113").Value = 6247.5
$ws.Range("K113").Value = 4499
$ws.Range("L113").Value = 6247.5
$ws.Range("M113").Value = -2329
$ws.Range("N113").Value = -10587.5

$ws.Range("H122").Value = 11286.917
$ws.Range("I122").Value = 4949.1665
$ws.Range("K122").Value = 14847.4995
$ws.Range("M122").Value = -12397.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2610686.8
$ws.Range("J81").Value = 15995
$ws.Range("L81").Value = 31990
$ws.Range("N81").Value = -34112

$ws.Range("H84").Value = 2610686.8
$ws.Range("J84").Value = 15995
$ws.Range("L84").Value = 159950
$ws.Range("N84").Value = -170558

$ws.Range("H132").Value = 26320708
$ws.Range("I132").Value = 1319.7059
$ws.Range("J132").Value = 250035500
$ws.Range("K132").Value = 3959.1177
$ws.Range("L132").Value = 750106500
$ws.Range("M132").Value = -1429.1177
$ws.Range("N132").Value = -750111560
